$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.459.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.69%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.430.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.46%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +1.33%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'308.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.26%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'88.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -8.14%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.41%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +1.15%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -7.46%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'31.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -8.49%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.29%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.41%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'2.813.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.21%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.37%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.424.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.43%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'14.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.77%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.753"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.22%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'41.067.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.82%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.50%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0890"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.06%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'68.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.07%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'10.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -10.77%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'230.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.75%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.44%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +0.14%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.85%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'23.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.09%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -1.23%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'9.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.80%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'34.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -7.69%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'151.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.67%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'5.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.09%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -3.39%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -4.73%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0729"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.34%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'16.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.67%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.90%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -8.71%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.50%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.0974"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -8.72%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.48%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.19%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'20.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.82%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.915.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.69%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.98%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -9.79%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.692.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.08%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'8.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.99%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'92.99"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'0.169"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -7.98%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'70.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -9.58%  "
$ws.Range("E51").Style = "Normal"
